# Refresh the cryptocurrency table's "Price" (D) and "Volume(1h)" (E) columns
# with the latest scrape, matching the GitHub Actions commit that updated
# cryptos.xlsx. Values are written with a leading apostrophe so Excel keeps
# them as literal text (the source column stores numbers/percentages as
# strings, e.g. "28.174.49", "0.5133", "  +0.39%  "), then the style is
# reset back to Normal so no stray quote-prefix formatting sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '28.174.49' }
    @{ Cell = "E2"; Value = '  +0.39%  ' }
    @{ Cell = "D3"; Value = '1.880.15' }
    @{ Cell = "E3"; Value = '  +0.39%  ' }
    @{ Cell = "E4"; Value = '  -0.14%  ' }
    @{ Cell = "E5"; Value = '  +0.08%  ' }
    @{ Cell = "E6"; Value = '  -0.13%  ' }
    @{ Cell = "D7"; Value = '0.5133' }
    @{ Cell = "E7"; Value = '  -0.06%  ' }
    @{ Cell = "D8"; Value = '0.3904' }
    @{ Cell = "E8"; Value = '  +1.99%  ' }
    @{ Cell = "D9"; Value = '0.08374' }
    @{ Cell = "E9"; Value = '  +1.38%  ' }
    @{ Cell = "D10"; Value = '1.116' }
    @{ Cell = "E10"; Value = '  +0.21%  ' }
    @{ Cell = "D11"; Value = '6.232' }
    @{ Cell = "E11"; Value = '  +0.15%  ' }
    @{ Cell = "E12"; Value = '  +1.19%  ' }
    @{ Cell = "D13"; Value = '1.873.46' }
    @{ Cell = "E13"; Value = '  -0.20%  ' }
    @{ Cell = "D14"; Value = '7.297' }
    @{ Cell = "E14"; Value = '  +1.25%  ' }
    @{ Cell = "E15"; Value = '  -0.18%  ' }
    @{ Cell = "D16"; Value = '0.00001107' }
    @{ Cell = "E16"; Value = '  +1.21%  ' }
    @{ Cell = "D17"; Value = '91.25' }
    @{ Cell = "E17"; Value = '  +0.20%  ' }
    @{ Cell = "D18"; Value = '0.06658' }
    @{ Cell = "E18"; Value = '  +0.10%  ' }
    @{ Cell = "D19"; Value = '17.73' }
    @{ Cell = "E19"; Value = '  -1.58%  ' }
    @{ Cell = "D20"; Value = '1.002' }
    @{ Cell = "E20"; Value = '  -0.10%  ' }
    @{ Cell = "D21"; Value = '6.052' }
    @{ Cell = "D22"; Value = '28.221.10' }
    @{ Cell = "D23"; Value = '11.19' }
    @{ Cell = "E23"; Value = '  +0.15%  ' }
    @{ Cell = "D24"; Value = '2.261' }
    @{ Cell = "E24"; Value = '  -0.15%  ' }
    @{ Cell = "D25"; Value = '2.088.05' }
    @{ Cell = "E25"; Value = '  -0.29%  ' }
    @{ Cell = "D26"; Value = '2.502' }
    @{ Cell = "E26"; Value = '  -3.33%  ' }
    @{ Cell = "D27"; Value = '158.41' }
    @{ Cell = "E27"; Value = '  +0.59%  ' }
    @{ Cell = "D28"; Value = '20.62' }
    @{ Cell = "E28"; Value = '  +0.10%  ' }
    @{ Cell = "D29"; Value = '125.26' }
    @{ Cell = "E29"; Value = '  -0.47%  ' }
    @{ Cell = "D30"; Value = '0.1065' }
    @{ Cell = "E30"; Value = '  +0.45%  ' }
    @{ Cell = "E31"; Value = '  -0.47%  ' }
    @{ Cell = "D32"; Value = '5.899' }
    @{ Cell = "E32"; Value = '  +5.20%  ' }
    @{ Cell = "D33"; Value = '3.591' }
    @{ Cell = "E33"; Value = '  -0.65%  ' }
    @{ Cell = "D34"; Value = '9.744' }
    @{ Cell = "E34"; Value = '  +0.74%  ' }
    @{ Cell = "E35"; Value = '  -0.08%  ' }
    @{ Cell = "D36"; Value = '0.06554' }
    @{ Cell = "E36"; Value = '  -0.65%  ' }
    @{ Cell = "D37"; Value = '0.2196' }
    @{ Cell = "E37"; Value = '  +1.40%  ' }
    @{ Cell = "D38"; Value = '1.210' }
    @{ Cell = "E38"; Value = '  -0.69%  ' }
    @{ Cell = "D39"; Value = '0.6528' }
    @{ Cell = "E39"; Value = '  +0.93%  ' }
    @{ Cell = "D40"; Value = '5.022' }
    @{ Cell = "E40"; Value = '  +2.88%  ' }
    @{ Cell = "D41"; Value = '1.230' }
    @{ Cell = "E41"; Value = '  -1.57%  ' }
    @{ Cell = "D42"; Value = '11.32' }
    @{ Cell = "E42"; Value = '  -0.24%  ' }
    @{ Cell = "D43"; Value = '0.6123' }
    @{ Cell = "E43"; Value = '  +0.22%  ' }
    @{ Cell = "D44"; Value = '13.12' }
    @{ Cell = "E44"; Value = '  +0.85%  ' }
    @{ Cell = "E45"; Value = '  -0.68%  ' }
    @{ Cell = "D46"; Value = '3.669' }
    @{ Cell = "E46"; Value = '  -0.16%  ' }
    @{ Cell = "D47"; Value = '2.015' }
    @{ Cell = "E47"; Value = '  +0.24%  ' }
    @{ Cell = "D48"; Value = '1.227' }
    @{ Cell = "E48"; Value = '  +0.64%  ' }
    @{ Cell = "D49"; Value = '121.59' }
    @{ Cell = "E49"; Value = '  +0.55%  ' }
    @{ Cell = "D50"; Value = '78.19' }
    @{ Cell = "E50"; Value = '  -2.95%  ' }
    @{ Cell = "D51"; Value = '0.06903' }
    @{ Cell = "E51"; Value = '  +0.21%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "`'$($u.Value)"
    $range.Style = "Normal"
}
